# "bugs, referral system - done"
#
# Fix the countryId value on the (only) data row of Sheet1: the account
# row was tagged for the USA but should be tagged for Germany (DE).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# A1:L1 holds the column headers (countryId, login, password, ...); A2:H2
# holds the single data row. Correct the country code bug.
$ws.Range("A2").Value = "DE"

# Leave the cursor on the cell that was just fixed, matching where the
# author's selection ended up after making the change.
$ws.Range("A2").Select()
